# Scheduled runner update: refresh cached market-board price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a batch of leve rows
# across the crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 40000
$ws.Range("I68").Value = 30000
$ws.Range("J68").Value = 50000
$ws.Range("K68").Value = 30000
$ws.Range("L68").Value = 50000
$ws.Range("M68").Value = -29251
$ws.Range("N68").Value = -51498

$ws.Range("H71").Value = 40000
$ws.Range("I71").Value = 30000
$ws.Range("J71").Value = 50000
$ws.Range("K71").Value = 90000
$ws.Range("L71").Value = 150000
$ws.Range("M71").Value = -86256
$ws.Range("N71").Value = -157488

$ws.Range("H86").Value = 127927.875
$ws.Range("I86").Value = 334691.34
$ws.Range("J86").Value = 3869.8
$ws.Range("K86").Value = 334691.34
$ws.Range("L86").Value = 3869.8
$ws.Range("M86").Value = -333568.34
$ws.Range("N86").Value = -6115.8

$ws.Range("H88").Value = 6081.2856
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H89").Value = 127927.875
$ws.Range("I89").Value = 334691.34
$ws.Range("J89").Value = 3869.8
$ws.Range("K89").Value = 1673456.7
$ws.Range("L89").Value = 19349
$ws.Range("M89").Value = -1667840.7
$ws.Range("N89").Value = -30581

$ws.Range("H91").Value = 6081.2856
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H92").Value = 46193.547
$ws.Range("I92").Value = 83997.086
$ws.Range("J92").Value = 829.3
$ws.Range("K92").Value = 83997.086
$ws.Range("L92").Value = 829.3
$ws.Range("M92").Value = -82749.086
$ws.Range("N92").Value = -3325.3

$ws.Range("H96").Value = 1388
$ws.Range("J96").Value = 2600
$ws.Range("L96").Value = 7800
$ws.Range("N96").Value = -10546

$ws.Range("H97").Value = 2715
$ws.Range("J97").Value = 3730.8333
$ws.Range("L97").Value = 11192.4999
$ws.Range("N97").Value = -12184.4999

$ws.Range("H98").Value = 1027.3889
$ws.Range("I98").Value = 1023.1177
$ws.Range("K98").Value = 1023.1177
$ws.Range("M98").Value = 474.8823

$ws.Range("H99").Value = 2304.0908
$ws.Range("J99").Value = 5969.25
$ws.Range("L99").Value = 17907.75
$ws.Range("N99").Value = -20903.75

$ws.Range("H100").Value = 2835.04
$ws.Range("I100").Value = 2127
$ws.Range("K100").Value = 2127
$ws.Range("M100").Value = -1586

$ws.Range("H122").Value = 1027.3889
$ws.Range("I122").Value = 1023.1177
$ws.Range("K122").Value = 3069.3531
$ws.Range("M122").Value = -619.3531000000003

$ws.Range("H138").Value = 3373.5833
$ws.Range("I138").Value = 2351.037
$ws.Range("K138").Value = 7053.110999999999
$ws.Range("M138").Value = -1913.110999999999

$ws.Range("H141").Value = 4680.5625
$ws.Range("I141").Value = 4680.5625
$ws.Range("K141").Value = 14041.6875
$ws.Range("M141").Value = -8861.6875


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3947.5
$ws.Range("I45").Value = 3171.5833
$ws.Range("K45").Value = 3171.5833
$ws.Range("M45").Value = -2794.5833

$ws.Range("H97").Value = 6535.7827
$ws.Range("I97").Value = 6612.4443
$ws.Range("J97").Value = 6259.8
$ws.Range("K97").Value = 6612.4443
$ws.Range("L97").Value = 6259.8
$ws.Range("M97").Value = -6116.4443
$ws.Range("N97").Value = -7251.8


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3148.4285
$ws.Range("I20").Value = 2875.4443
$ws.Range("J20").Value = 3639.8
$ws.Range("K20").Value = 2875.4443
$ws.Range("L20").Value = 3639.8
$ws.Range("M20").Value = -2628.4443
$ws.Range("N20").Value = -4133.8

$ws.Range("H86").Value = 4056.2856
$ws.Range("I86").Value = 3678.8
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3678.8
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2555.8
$ws.Range("N86").Value = -7246

$ws.Range("H89").Value = 4056.2856
$ws.Range("I89").Value = 3678.8
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 18394
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -12778
$ws.Range("N89").Value = -36232

$ws.Range("H94").Value = 3452.7693
$ws.Range("I94").Value = 2899.9
$ws.Range("J94").Value = 5295.6665
$ws.Range("K94").Value = 2899.9
$ws.Range("L94").Value = 5295.6665
$ws.Range("M94").Value = -2448.9
$ws.Range("N94").Value = -6197.6665

$ws.Range("H99").Value = 30399.473
$ws.Range("I99").Value = 37788.25
$ws.Range("K99").Value = 37788.25
$ws.Range("M99").Value = -36290.25

$ws.Range("H107").Value = 2300
$ws.Range("I107").Value = 1400
$ws.Range("K107").Value = 1400
$ws.Range("M107").Value = 520


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 524.17645
$ws.Range("I107").Value = 517
$ws.Range("J107").Value = 547.5
$ws.Range("K107").Value = 517
$ws.Range("L107").Value = 547.5
$ws.Range("M107").Value = 1403
$ws.Range("N107").Value = -4387.5

$ws.Range("H134").Value = 41554.5
$ws.Range("I134").Value = 54774
$ws.Range("K134").Value = 164322
$ws.Range("M134").Value = -161787


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 502
$ws.Range("J12").Value = 198.7
$ws.Range("L12").Value = 596.0999999999999
$ws.Range("N12").Value = -942.0999999999999

$ws.Range("H25").Value = 6698.4
$ws.Range("I25").Value = 1245
$ws.Range("J25").Value = 10334
$ws.Range("K25").Value = 3735
$ws.Range("L25").Value = 31002
$ws.Range("M25").Value = -3566
$ws.Range("N25").Value = -31340

$ws.Range("H30").Value = 6698.4
$ws.Range("I30").Value = 1245
$ws.Range("J30").Value = 10334
$ws.Range("K30").Value = 3735
$ws.Range("L30").Value = 31002
$ws.Range("M30").Value = -3633
$ws.Range("N30").Value = -31206

$ws.Range("H64").Value = 5001.5713
$ws.Range("J64").Value = 7252.75
$ws.Range("L64").Value = 21758.25
$ws.Range("N64").Value = -22298.25

$ws.Range("H67").Value = 5001.5713
$ws.Range("J67").Value = 7252.75
$ws.Range("L67").Value = 21758.25
$ws.Range("N67").Value = -23630.25

$ws.Range("H113").Value = 967.6875
$ws.Range("I113").Value = 1344.5
$ws.Range("K113").Value = 4033.5
$ws.Range("M113").Value = -1863.5

$ws.Range("H128").Value = 131007
$ws.Range("I128").Value = 131007
$ws.Range("K128").Value = 393021
$ws.Range("M128").Value = -388041


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9241.916999999999
$ws.Range("I70").Value = 8211.556
$ws.Range("K70").Value = 8211.556
$ws.Range("M70").Value = -7941.556

$ws.Range("H73").Value = 9241.916999999999
$ws.Range("I73").Value = 8211.556
$ws.Range("K73").Value = 8211.556
$ws.Range("M73").Value = -7275.556

$ws.Range("H80").Value = 2748.72
$ws.Range("I80").Value = 2310.4546
$ws.Range("J80").Value = 3093.0715
$ws.Range("K80").Value = 2310.4546
$ws.Range("L80").Value = 3093.0715
$ws.Range("M80").Value = -1312.4546
$ws.Range("N80").Value = -5089.0715

$ws.Range("H83").Value = 2748.72
$ws.Range("I83").Value = 2310.4546
$ws.Range("J83").Value = 3093.0715
$ws.Range("K83").Value = 11552.273
$ws.Range("L83").Value = 15465.3575
$ws.Range("M83").Value = -6560.273000000001
$ws.Range("N83").Value = -25449.3575

$ws.Range("H132").Value = 32737
$ws.Range("I132").Value = 40818.69
$ws.Range("J132").Value = 2719.2856
$ws.Range("K132").Value = 122456.07
$ws.Range("L132").Value = 8157.8568
$ws.Range("M132").Value = -119926.07
$ws.Range("N132").Value = -13217.8568


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9485.913
$ws.Range("I7").Value = 15269.3
$ws.Range("J7").Value = 5037.154
$ws.Range("K7").Value = 15269.3
$ws.Range("L7").Value = 5037.154
$ws.Range("M7").Value = -15157.3
$ws.Range("N7").Value = -5261.154

$ws.Range("H126").Value = 9485.913
$ws.Range("I126").Value = 15269.3
$ws.Range("J126").Value = 5037.154
$ws.Range("K126").Value = 45807.89999999999
$ws.Range("L126").Value = 15111.462
$ws.Range("M126").Value = -43337.89999999999
$ws.Range("N126").Value = -20051.462

$ws.Range("H132").Value = 45901
$ws.Range("I132").Value = 54661.914
$ws.Range("K132").Value = 163985.742
$ws.Range("M132").Value = -161455.742


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 36662.668
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 36662.668
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 36662.668
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -38160.668

$ws.Range("H72").Value = 36662.668
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 36662.668
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 109988.004
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -117476.004

$ws.Range("H123").Value = 99998
$ws.Range("J123").Value = 99998
$ws.Range("L123").Value = 99998
$ws.Range("N123").Value = -109798
